$d = $word.ActiveDocument

# Locate the empty "True or False" list paragraph that currently holds
# only the _GoBack bookmark (the last list item under numId=2).
$count = $d.Paragraphs.Count
$targetIndex = -1
for ($i = 1; $i -le $count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text.Length -le 1 -and $p.Range.ListFormat.ListType -ne 0) {
        $targetIndex = $i
    }
}

if ($targetIndex -eq -1) {
    throw "Could not locate target paragraph"
}

$target = $d.Paragraphs.Item($targetIndex)
$r = $target.Range
$r.Collapse(1)

$xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">In C#, you must declare a variable in a program before you can use it to store data. </w:t></w:r><w:r><w:rPr><w:b/></w:rPr><w:t xml:space="preserve">Answer: </w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">You can declare multiple variables of different data types with one declaration. </w:t></w:r><w:r><w:rPr><w:b/></w:rPr><w:t xml:space="preserve">Answer: </w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">When you append the letter D or d to a numeric literal, it is treated as a decimal and is referred to as a decimal literal. </w:t></w:r><w:r><w:rPr><w:b/></w:rPr><w:t xml:space="preserve">Answer: </w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">The order of operations dictates that the division operator works before the addition operator does. </w:t></w:r><w:r><w:rPr><w:b/></w:rPr><w:t xml:space="preserve">Answer: </w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">All variables have a </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>ToString</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> method that you can call to convert the variable’s value to a string. </w:t></w:r><w:r><w:rPr><w:b/></w:rPr><w:t xml:space="preserve">Answer: </w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">When you pass the formatting string “C” or “c” to the </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>ToString</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> method, the number is returned formatted as currency. </w:t></w:r><w:r><w:rPr><w:b/></w:rPr><w:t xml:space="preserve">Answer: </w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">When you declare a named constant, an initialization value is required. </w:t></w:r><w:r><w:rPr><w:b/></w:rPr><w:t xml:space="preserve">Answer: </w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">An error will occur if the compiler finds a statement that tries to change the value of a constant field. </w:t></w:r><w:r><w:rPr><w:b/></w:rPr><w:t xml:space="preserve">Answer: </w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">Forms and most controls have a Preferences property that allows you to change the object’s background color. </w:t></w:r><w:r><w:rPr><w:b/></w:rPr><w:t xml:space="preserve">Answer: </w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$r.InsertXML($xml)
